$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.221053335166417
